# Grimstone Beasts - translation update
# Adds three new rows (139-141) to Sheet1 for the ThingDef "Grimstone_NorthernLynx"
# tool labels (left claw / right claw / head), and fixes two stray cell styles
# (F110, F133) that no longer need their one-off alignment-applied format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 139: tools.0.label (left claw) -----------------------------
$ws.Range("B139").Value = "ThingDef"
$ws.Range("F139").Value = "왼쪽 발톱"
$ws.Range("A139").Value = "ThingDef+Grimstone_NorthernLynx.tools.0.label"
$ws.Range("C139").Value = "Grimstone_NorthernLynx.tools.0.label"

# --- New row 140: tools.1.label (right claw) -----------------------------
$ws.Range("B140").Value = "ThingDef"
$ws.Range("F140").Value = "오른쪽 발톱"
$ws.Range("A140").Value = "ThingDef+Grimstone_NorthernLynx.tools.1.label"
$ws.Range("C140").Value = "Grimstone_NorthernLynx.tools.1.label"

# --- New row 141: tools.3.label (head) ------------------------------------
$ws.Range("B141").Value = "ThingDef"
$ws.Range("F141").Value = "머리"
$ws.Range("A141").Value = "ThingDef+Grimstone_NorthernLynx.tools.3.label"
$ws.Range("C141").Value = "Grimstone_NorthernLynx.tools.3.label"

# --- Tidy up two cells that previously carried a stray one-off style -----
# (F110 / F133 used a cellXf with applyAlignment="1" but no real alignment
# override; re-stamping the font name collapses them back onto the shared
# "Korean body text" style used by every other cell in the sheet.)
$ws.Range("F110").Font.Name = "맑은 고딕"
$ws.Range("F133").Font.Name = "맑은 고딕"

# --- Move the active selection the way the source workbook ended up ------
$ws.Range("L141").Select() | Out-Null
